$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded, so insert a row right
# before the current row 48 — this pushes the existing rows 48-53 down
# to 49-54 (matching the "shifted" rows in the diff) and grows the used
# range from A1:R53 to A1:R54.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new weekly record.
$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 44522
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = 100114002
$ws.Range("G48").Value = "Camote"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 1600
$ws.Range("K48").Value = 11000
$ws.Range("L48").Value = 12000
$ws.Range("M48").Value = 11500
$ws.Range("N48").Value = "$/malla 18 kilos"
$ws.Range("O48").Value = "Perú"
$ws.Range("P48").Value = 639
$ws.Range("Q48").Value = 18
$ws.Range("R48").Value = "Hortaliza"
